$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 2, pushing current row 2 (John Smith / us-001 invoice)
# down to row 3.
$ws.Rows.Item(2).Insert()

# Make sure the values that look like numbers/dates stay plain text (not
# auto-converted by Excel's type inference).
$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2:D2").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"

# Fill in the new row 2 with the new invoice data.
$ws.Range("A2").Value = "John Smith"
$ws.Range("B2").Value = "123456789"
$ws.Range("C2").Value = "2023-11-23"
$ws.Range("D2").Value = "2023-11-30"
$ws.Range("E2").Value = 40.25
$ws.Range("F2").Value = 2.05
$ws.Range("G2").Value = 42.3
$ws.Range("H2").Value = "Unpaid"

# Correct the invoice date on the row that shifted down to row 3.
$ws.Range("C3").Value = "2019-11-02"
